$wb = $excel.ActiveWorkbook

# --- DBS sheet: add a new "txtNoLastFirst" row under the existing data ---
$ws2 = $wb.Worksheets.Item("DBS")
$ws2.Range("A3").Value = "txtNoLastFirst"
$ws2.Range("B3").Value = "TlrNo = "
$ws2.Range("C3").Value = "CreateDate desc"

# --- DBD sheet: move the cursor/selection (view-state only change) ---
$ws1 = $wb.Worksheets.Item("DBD")
$ws1.Activate() | Out-Null
$ws1.Range("B18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# --- DBS sheet: re-activate and move the cursor/selection to the new row ---
$ws2.Activate() | Out-Null
$ws2.Range("B3").Select() | Out-Null
